$wb = $excel.ActiveWorkbook

# 1) Rename first sheet: "Vista02" -> "Alojamiento"
$wsAlojamiento = $wb.Worksheets.Item(1)
$wsAlojamiento.Name = "Alojamiento"

# 2) Vista04 (3rd sheet): selection moves from C1 to whole column B (active cell B1)
$wsVista04 = $wb.Worksheets.Item(3)
$wsVista04.Range("B1:B1048576").Select()

# 3) Vista05 (4th sheet): delete the leading (mostly blank) column A so B/C/D shift to A/B/C,
#    restore the "Año" header label in the now-first column, and become the active/selected tab.
$wsVista05 = $wb.Worksheets.Item(4)
$wsVista05.Range("A1").EntireColumn.Delete()
$wsVista05.Range("A1").Value = "Año"
$wsVista05.Activate()
$wsVista05.Range("F18").Select()
